$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the "R10" rule row from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection on the sheet at the edited cell
$ws.Range("E8").Select()
